# Update "想去人数" (F column) values on the "展览" (Exhibition) and
# "全部类型" (All Types) sheets to match the newly generated data snapshot.
# Same set of events are updated on both sheets, but the events sit on
# different row numbers on each sheet, so two row maps are used.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$sheet1Updates = @{
    6  = 14355
    7  = 16603
    9  = 103
    10 = 15
    12 = 203
    23 = 70
    27 = 6753
    28 = 973
    31 = 12
    33 = 5758
    35 = 145
    37 = 4836
}

# Sheet "全部类型": row -> new F value
$sheet4Updates = @{
    6  = 14355
    7  = 16603
    9  = 103
    10 = 15
    12 = 203
    23 = 70
    28 = 6753
    29 = 973
    32 = 12
    36 = 5758
    38 = 145
    40 = 4836
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
